$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row reference (column A) swaps caused by re-sorting a handful of countries
# within the shared-string/country list, plus the refreshed COVID-19
# numeric stats (columns B-H) for the affected rows, and the
# "Datos actualizados" timestamp in A1.
$changes = @(
    @{ Cell = "A48"; Value = "Barein" },
    @{ Cell = "A49"; Value = "Afganistan" },
    @{ Cell = "A50"; Value = "Rumania" },
    @{ Cell = "A69"; Value = "Costa de Marfil" },
    @{ Cell = "A70"; Value = "Chequia" },
    @{ Cell = "A110"; Value = "Malaui" },
    @{ Cell = "A111"; Value = "Sri Lanka" },
    @{ Cell = "A128"; Value = "Yemen" },
    @{ Cell = "A129"; Value = "Nueva Zelanda" },
    @{ Cell = "A147"; Value = "Surinam" },
    @{ Cell = "A148"; Value = "Republica del Chad" },
    @{ Cell = "A149"; Value = "Principado de Andorra" },
    @{ Cell = "A209"; Value = "Groenlandia" },
    @{ Cell = "A210"; Value = "Islas Malvinas" },
    @{ Cell = "B4"; Value = 3677325 },
    @{ Cell = "C4"; Value = 60498 },
    @{ Cell = "D4"; Value = 1671924 },
    @{ Cell = "E4"; Value = 1864481 },
    @{ Cell = "G4"; Value = 776 },
    @{ Cell = "H4"; Value = 140920 },
    @{ Cell = "B5"; Value = 2012151 },
    @{ Cell = "C5"; Value = 41242 },
    @{ Cell = "E5"; Value = 568688 },
    @{ Cell = "G5"; Value = 1165 },
    @{ Cell = "H5"; Value = 76688 },
    @{ Cell = "B19"; Value = 201832 },
    @{ Cell = "C19"; Value = 580 },
    @{ Cell = "E19"; Value = 6275 },
    @{ Cell = "B27"; Value = 85771 },
    @{ Cell = "C27"; Value = 928 },
    @{ Cell = "D27"; Value = 26691 },
    @{ Cell = "E27"; Value = 54960 },
    @{ Cell = "G27"; Value = 53 },
    @{ Cell = "H27"; Value = 4120 },
    @{ Cell = "B46"; Value = 46059 },
    @{ Cell = "C46"; Value = 1871 },
    @{ Cell = "D46"; Value = 20370 },
    @{ Cell = "E46"; Value = 25305 },
    @{ Cell = "G46"; Value = 8 },
    @{ Cell = "H46"; Value = 384 },
    @{ Cell = "B48"; Value = 35084 },
    @{ Cell = "C48"; Value = 524 },
    @{ Cell = "D48"; Value = 30809 },
    @{ Cell = "E48"; Value = 4154 },
    @{ Cell = "G48"; Value = 4 },
    @{ Cell = "H48"; Value = 121 },
    @{ Cell = "B49"; Value = 35070 },
    @{ Cell = "C49"; Value = 76 },
    @{ Cell = "D49"; Value = 22824 },
    @{ Cell = "E49"; Value = 11133 },
    @{ Cell = "H49"; Value = 1113 },
    @{ Cell = "B50"; Value = 35003 },
    @{ Cell = "C50"; Value = 777 },
    @{ Cell = "D50"; Value = 22189 },
    @{ Cell = "E50"; Value = 10843 },
    @{ Cell = "G50"; Value = 19 },
    @{ Cell = "H50"; Value = 1971 },
    @{ Cell = "D53"; Value = 29900 },
    @{ Cell = "E53"; Value = 1421 },
    @{ Cell = "B69"; Value = 13554 },
    @{ Cell = "C69"; Value = 151 },
    @{ Cell = "D69"; Value = 7363 },
    @{ Cell = "E69"; Value = 6104 },
    @{ Cell = "H69"; Value = 87 },
    @{ Cell = "B70"; Value = 13551 },
    @{ Cell = "C70"; Value = 76 },
    @{ Cell = "D70"; Value = 8629 },
    @{ Cell = "E70"; Value = 4567 },
    @{ Cell = "H70"; Value = 355 },
    @{ Cell = "B96"; Value = 5285 },
    @{ Cell = "C96"; Value = 163 },
    @{ Cell = "D96"; Value = 4275 },
    @{ Cell = "E96"; Value = 899 },
    @{ Cell = "B110"; Value = 2712 },
    @{ Cell = "C110"; Value = 98 },
    @{ Cell = "D110"; Value = 1073 },
    @{ Cell = "E110"; Value = 1588 },
    @{ Cell = "G110"; Value = 8 },
    @{ Cell = "H110"; Value = 51 },
    @{ Cell = "B111"; Value = 2686 },
    @{ Cell = "C111"; Value = 15 },
    @{ Cell = "D111"; Value = 2007 },
    @{ Cell = "E111"; Value = 668 },
    @{ Cell = "H111"; Value = 11 },
    @{ Cell = "B128"; Value = 1552 },
    @{ Cell = "C128"; Value = 26 },
    @{ Cell = "D128"; Value = 695 },
    @{ Cell = "E128"; Value = 419 },
    @{ Cell = "G128"; Value = 5 },
    @{ Cell = "H128"; Value = 438 },
    @{ Cell = "B129"; Value = 1548 },
    @{ Cell = "C129"; Value = 1 },
    @{ Cell = "D129"; Value = 1499 },
    @{ Cell = "E129"; Value = 27 },
    @{ Cell = "H129"; Value = 22 },
    @{ Cell = "B131"; Value = 1473 },
    @{ Cell = "C131"; Value = 38 },
    @{ Cell = "D131"; Value = 770 },
    @{ Cell = "E131"; Value = 699 },
    @{ Cell = "B147"; Value = 886 },
    @{ Cell = "C147"; Value = 49 },
    @{ Cell = "D147"; Value = 581 },
    @{ Cell = "E147"; Value = 287 },
    @{ Cell = "H147"; Value = 18 },
    @{ Cell = "B148"; Value = 886 },
    @{ Cell = "C148"; Value = 1 },
    @{ Cell = "D148"; Value = 799 },
    @{ Cell = "E148"; Value = 12 },
    @{ Cell = "H148"; Value = 75 },
    @{ Cell = "B149"; Value = 862 },
    @{ Cell = "C149"; Value = 0 },
    @{ Cell = "D149"; Value = 803 },
    @{ Cell = "E149"; Value = 7 },
    @{ Cell = "H149"; Value = 52 },
    @{ Cell = "B157"; Value = 607 },
    @{ Cell = "C157"; Value = 31 },
    @{ Cell = "E157"; Value = 455 },
    @{ Cell = "G157"; Value = 1 },
    @{ Cell = "H157"; Value = 28 }
)

foreach ($ch in $changes) {
    $ws.Range($ch.Cell).Value = $ch.Value
}

$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 23:31"

Write-Host ("Applied " + $changes.Count + " cell updates plus timestamp.")
